# Adds a new slide (38th) to the deck: a title textbox + a 4-column,
# 21-row (1 header + 20 data) table of top DiSCoVER drugs for the
# cerebellar stem cell control comparison.

$p = $ppt.ActivePresentation

# New slide goes at the end (index 38), using the same blank layout
# (slideLayout7) that every other slide in this deck uses.
$s = $p.Slides.Add($p.Slides.Count + 1, 7)

# ---- Title textbox -------------------------------------------------
$titleText = 'DiSCoVER: top drugs (cerebellar stem cell control)'
$tb = $s.Shapes.AddTextbox(1, 0, 0, 720, 54)
$tb.Name = 'TextBox 1'
$tb.Fill.Visible = $false
$tb.TextFrame.WordWrap = $false
$tf = $tb.TextFrame.TextRange
$tf.Text = "`r" + $titleText
$tf.Paragraphs(2).Font.Size = 26
$tb.TextFrame.AutoSize = 1
$tb.Height = 54

# ---- Drug table ------------------------------------------------------
$data = @(
    @('tl-2-105', '0.65', '+..', 'Not Clinically Relevant'),
    @('sb52334', '0.62', '+..', 'Not Clinically Relevant'),
    @('gsk1070916', '0.61', '+..', 'Not Clinically Relevant'),
    @('ql-xii-61', '0.59', '+..', 'Not Clinically Relevant'),
    @('linsitinib', '0.55', '++.', 'IGF-1R inhibitor'),
    @('gw-2580', '0.55', '+..', 'Not Clinically Relevant'),
    @('gsk429286a', '0.55', '+..', 'Not Clinically Relevant'),
    @('tubastatin a', '0.54', '++.', 'Not Clinically Relevant'),
    @('bx-912', '0.54', '+..', 'Not Clinically Relevant'),
    @('vx-702', '0.53', '+..', 'Not Clinically Relevant'),
    @('rucaparib', '0.52', '+..', 'PARP inhibitor, inhibits DNA repair'),
    @('navitoclax', '0.51', '++.', 'Bcl-2 family inhibitor: esp Bcl-xL, Bcl-2 and Bcl-w'),
    @('xmd14-99', '0.48', '+..', 'Not Clinically Relevant'),
    @('axitinib', '0.47', '++.', 'VEGFR, c-KIT and PDGFR inhibitor'),
    @('hg-5-88-01', '0.47', '+..', 'Not Clinically Relevant'),
    @('amuvatinib', '0.46', '+..', 'Not Clinically Relevant'),
    @('xmd13-2', '0.45', '+..', 'Not Clinically Relevant'),
    @('indisulam', '0.44', '.+.', 'carbonic anhydrase inibitor and  CDK inhibitor, targets G1 by depleting cyclin E. inducing p53 and p21, and inhibiting CDK2'),
    @('nsc-87877', '0.43', '+..', 'Not Clinically Relevant'),
    @('cil55', '0.43', '.+.', 'Not Clinically Relevant')
)

$headers = @('Drug', 'Score', 'Evidence', 'Mechanism of action')

$nRows = $data.Count + 1
$tbl = $s.Shapes.AddTable($nRows, 4, 32.4, 61.2, 651.6, 324)
$tbl.Name = 'Table 2'

$colWidths = @(79.2, 57.6, 82.8, 432)
for ($c = 1; $c -le 4; $c++) {
    $tbl.Table.Columns.Item($c).Width = $colWidths[$c - 1]
}

for ($c = 1; $c -le 4; $c++) {
    $tbl.Table.Cell(1, $c).Shape.TextFrame.TextRange.Text = $headers[$c - 1]
}

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    for ($c = 1; $c -le 4; $c++) {
        $cellRange = $tbl.Table.Cell($r + 2, $c).Shape.TextFrame.TextRange
        $cellRange.Text = $row[$c - 1]
        $cellRange.Font.Size = 10.5
    }
}

# Nudge the final row's height a hair so the table's total height lands
# on exactly 4114800 EMU (matches the source row-height rounding).
$tbl.Table.Rows.Item($nRows).Height = 15.42992125984252

